# Updates the cryptocurrency price/volume table to reflect the latest
# GitHub Actions data refresh (source data keeps every Price-column cell
# as plain text, even the ones that happen to look numeric, e.g. "515.28").
# Excel's COM layer auto-coerces a numeric-looking string into a real
# number on assignment, so for just those cells we flip NumberFormat to
# Text ("@") first to keep them as text like the rest of the column;
# values that can never parse as a number (e.g. "57.813.81") don't need
# that treatment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.813.81'
$ws.Range("E2").Value = '  -2.70%  '

# Row 3
$ws.Range("D3").Value = '2.566.25'
$ws.Range("E3").Value = '  -2.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.28'
$ws.Range("E5").Value = '  -2.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.22'
$ws.Range("E6").Value = '  -4.91%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.558'
$ws.Range("E8").Value = '  -2.26%  '

# Row 9
$ws.Range("D9").Value = '2.587.44'
$ws.Range("E9").Value = '  -2.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.42'
$ws.Range("E10").Value = '  -3.98%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0988'
$ws.Range("E11").Value = '  -5.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.324'
$ws.Range("E12").Value = '  -4.30%  '

# Row 13
$ws.Range("E13").Value = '  +0.34%  '

# Row 14
$ws.Range("D14").Value = '3.022.08'
$ws.Range("E14").Value = '  -2.87%  '

# Row 15
$ws.Range("D15").Value = '57.809.50'
$ws.Range("E15").Value = '  -2.62%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.97'
$ws.Range("E16").Value = '  -4.80%  '

# Row 17
$ws.Range("D17").Value = '2.570.00'
$ws.Range("E17").Value = '  -4.33%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000131'
$ws.Range("E18").Value = '  -4.62%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '331.00'
$ws.Range("E19").Value = '  -3.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.27'
$ws.Range("E20").Value = '  -4.53%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.03'
$ws.Range("E21").Value = '  -5.74%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.33'
$ws.Range("E22").Value = '  -1.19%  '

# Row 23
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.82'
$ws.Range("E24").Value = '  +0.18%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.165'
$ws.Range("E25").Value = '  -2.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.396'
$ws.Range("E27").Value = '  -5.55%  '

# Row 28
$ws.Range("D28").Value = '2.690.13'
$ws.Range("E28").Value = '  -2.48%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.88'
$ws.Range("E29").Value = '  -5.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.08%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0709'
$ws.Range("E31").Value = '  -11.79%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.89'
$ws.Range("E32").Value = '  -8.69%  '

# Row 33
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.55'
$ws.Range("E33").Value = '  -4.05%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.56'
$ws.Range("E34").Value = '  -2.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.02'
$ws.Range("E35").Value = '  -0.74%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.88'
$ws.Range("E36").Value = '  -7.77%  '

# Row 37
$ws.Range("E37").Value = '  -8.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.14'
$ws.Range("E38").Value = '  -1.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.829'
$ws.Range("E39").Value = '  -4.61%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("E40").Value = '  -8.16%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.41'
$ws.Range("E41").Value = '  -5.79%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.48'
$ws.Range("E42").Value = '  -4.99%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.07%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '271.55'
$ws.Range("E44").Value = '  +0.31%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.69'
$ws.Range("E45").Value = '  +0.31%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.588'
$ws.Range("E46").Value = '  -2.64%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0937'
$ws.Range("E47").Value = '  -3.94%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0512'
$ws.Range("E48").Value = '  -5.11%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.34'
$ws.Range("E49").Value = '  -5.74%  '

# Row 50
$ws.Range("D50").Value = '1.967.00'
$ws.Range("E50").Value = '  -3.49%  '

# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0216'
$ws.Range("E51").Value = '  -5.89%  '
